$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-11 form a block of records that has been rotated: each record has
# moved DOWN by one row (old row 2 -> row 3, old row 3 -> row 4, ...,
# old row 10 -> row 11), and the record that used to be in the last row (11)
# wraps around to become the first row (2).
#
# Only the columns whose values actually differ between records need to be
# touched (A, B, D, E, F, G, H, Q, R, S, Y, AA, AN, AO) - the remaining
# populated columns (C, I, P, T, U, V, W, Z, AB, AD, AE, AG, AJ, AK, AT, AW,
# AX, AY) hold the same value in every one of these rows, so rewriting them
# is a no-op that's better left alone.

$firstRow = 2
$lastRow = 11
$cols = @(1, 2, 4, 5, 6, 7, 8, 17, 18, 19, 25, 27, 40, 41)   # A,B,D,E,F,G,H,Q,R,S,Y,AA,AN,AO
$dateCols = @(25, 27)   # Y, AA hold plain-text dates ("2014-07-19"); writing them
                        # straight back through Value2 would get auto-parsed into a
                        # date serial, so round-trip them as explicit text instead.

# Stash the last row's values (the record that needs to wrap to the top).
$buffer = @{}
foreach ($c in $cols) {
    $buffer[$c] = $ws.Cells.Item($lastRow, $c).Value2
}

# Shift rows 2..10 down into rows 3..11 (process bottom-up so we never
# overwrite a source row before it has been read).
for ($r = $lastRow; $r -gt $firstRow; $r--) {
    $srcRow = $r - 1
    foreach ($c in $cols) {
        $v = $ws.Cells.Item($srcRow, $c).Value2
        $dst = $ws.Cells.Item($r, $c)
        $isDateLike = $dateCols -contains $c
        if ($isDateLike) {
            $dst.NumberFormat = "@"
            $dst.Value2 = $v
            $dst.Style = "Normal"
        } else {
            $dst.Value2 = $v
        }
    }
}

# Put the original last-row data into row 2.
foreach ($c in $cols) {
    $dst = $ws.Cells.Item($firstRow, $c)
    $isDateLike = $dateCols -contains $c
    if ($isDateLike) {
        $dst.NumberFormat = "@"
        $dst.Value2 = $buffer[$c]
        $dst.Style = "Normal"
    } else {
        $dst.Value2 = $buffer[$c]
    }
}
